$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 5
$ws.Range("C1").Value = 1623239774
$ws.Range("D1").Value = 72
$ws.Range("E1").Value = 1550

$ws.Range("A2").Value = 10
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 1623239781
$ws.Range("D2").Value = 89
$ws.Range("E2").Value = 7596

$ws.Range("I11").Select()
